$wb = $excel.ActiveWorkbook

# The workbook has duplicate data in "展览" and "全部类型" sheets; update both.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 441
    $ws.Range("F4").Value = 3213
    $ws.Range("F6").Value = 644
}
